$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A so existing data shifts right (B:E -> C:F)
$ws.Range("A1:A1").EntireColumn.Insert()

# New header for column A
$ws.Range("A1").Value = "Prediction Horizon"
$ws.Range("B1").Copy()
$ws.Range("A1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Insert 3 new rows above current row 2 (after the header) for horizon = 1 data
$ws.Range("A2:A4").EntireRow.Insert()
$ws.Range("A2:F4").ClearFormats()

# Fill in horizon = 1 rows (new rows 2-4)
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "persistence"
$ws.Range("C2").Value = 0.02046225165562914
$ws.Range("D2").Value = 0.0008541031879424526
$ws.Range("E2").Value = 0.02922504384842652
$ws.Range("F2").Value = 0.9832394036160876

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "random_forest"
$ws.Range("C3").Value = 0.005373401621674758
$ws.Range("D3").Value = 0.0001075083204025303
$ws.Range("E3").Value = 0.01036862191433993
$ws.Range("F3").Value = 0.9978876762261663

$ws.Range("A4").Value = 1
$ws.Range("B4").Value = "neural_network"
$ws.Range("C4").Value = 0.009708198715866524
$ws.Range("D4").Value = 0.0002432799561024466
$ws.Range("E4").Value = 0.01559743427947195
$ws.Range("F4").Value = 0.9950598850805302

# Fill in the Prediction Horizon (column A) for the existing horizon = 6 rows (now rows 5-7)
$ws.Range("A5").Value = 6
$ws.Range("A6").Value = 6
$ws.Range("A7").Value = 6

# Update the neural_network row (row 7) error values per diff
$ws.Range("C7").Value = 0.06599178062754969
$ws.Range("D7").Value = 0.009445472209663417
$ws.Range("E7").Value = 0.09718781924533247
$ws.Range("F7").Value = 0.776892828741772
